# Applies the "finalize camera ready version" edits described by the diff.
# All four edits are pure text-content changes inside the "Response:" runs;
# since the runs being merged/split share identical run formatting, a plain
# Find & Replace (literal text, no wildcards) reproduces the visible content
# change without altering formatting.

$d = $word.ActiveDocument

function Replace-Literal($find, $replace) {
    $r = $d.Content
    $ok = $r.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $find"
    }
}

# 1) "In the revised manuscript, " / "the " / "introduction is re-written..."
#    runs merge into a single run with identical text (no wording change).
Replace-Literal `
    "In the revised manuscript, the introduction is re-written to focus on the big picture of MI segmentation and how the proposed methods can address the current limitations." `
    "In the revised manuscript, the introduction is re-written to focus on the big picture of MI segmentation and how the proposed methods can address the current limitations."

# 2) "...the algorithm perform" / "s" / " better. ..." runs merge into a
#    single run with identical text (no wording change).
Replace-Literal `
    "the algorithm performs better. But multiple" `
    "the algorithm performs better. But multiple"

# 3) Replace the latter half of the AUPR paragraph with new reviewer-response text.
Replace-Literal `
    "However, for BRNN and CRNN, the POS tags also lowered the AUPR values.  In the 2nd paragraph of the discussion section, a discussion was given to explain why  MLP has the best precision and POS tags lowered the AUPR values of BRNN and CRNN. " `
    ("We believe that the BRNN and CRNN performed poorly with POS features because POS tagging is a supervised learning solution that uses features like the previous and next word. Since we already considered neighbor words by utilizing bi-directional RNN, it failed to achieve good results with redundant information. We observed that MLP achieved the highest precision which may be related to the fact that MLP poorly learned " + [char]0x201C + "new segment" + [char]0x201D + " and misclassified new segment words to same segments in 30%-40% of the time. ")

# 4) Replace "The caption font was modified" with "We modified the figure 3".
Replace-Literal `
    "The caption font was modified in the revised manuscript." `
    "We modified the figure 3 in the revised manuscript."
